# Generate Report for Handback
# Fills in the "Latest Target File", "Latest Handback File" and
# "Latest Handback DateTime" columns for the zh-cn and de-de handback
# tables now that both locales came back in sync with en-US, and flips
# the Status column (and the Overview rollup columns) from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/27257558a935de3c9d0b4a666fb36a88cdc689be/e2e/"

$statusHandedBack = "Handed back: in sync with en-US"

$file1 = "972394dd-bbf8-458c-a120-f79373f3dcf2"
$file2 = "aa87a066-626a-4b3f-af63-a554a7d50982"

$file1Md  = "$file1.md"
$file2Md  = "$file2.md"

# ---------------------------------------------------------------
# Overview sheet: status rollup columns for zh-cn (E) and de-de (F)
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C3").Value = $statusHandedBack

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($repoBase + $file1Md), "", "", $file1Md)
$zhcn.Range("J2").Value = "$file1.c2e4f7dcfa672628916e5aac907833914a96acc6.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-03 04:33:08"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($repoBase + $file2Md), "", "", $file2Md)
$zhcn.Range("J3").Value = "$file2.e5cf53f6cd40c754d025aeeab221fedbffa1ec43.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-03 04:33:08"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C3").Value = $statusHandedBack

$dede.Hyperlinks.Add($dede.Range("I2"), ($repoBase + $file1Md), "", "", $file1Md)
$dede.Range("J2").Value = "$file1.c2e4f7dcfa672628916e5aac907833914a96acc6.de-de.xlf"
$dede.Range("K2").Value = "2016-09-03 04:33:15"

$dede.Hyperlinks.Add($dede.Range("I3"), ($repoBase + $file2Md), "", "", $file2Md)
$dede.Range("J3").Value = "$file2.e5cf53f6cd40c754d025aeeab221fedbffa1ec43.de-de.xlf"
$dede.Range("K3").Value = "2016-09-03 04:33:15"

# ---------------------------------------------------------------
# Widen the columns that now hold the longer status / file-name text
# ---------------------------------------------------------------
$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40
